$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "05/25/2025"
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = 0.0004602700000000005
$ws.Range("C18").Value = 108631.8899776217
$ws.Range("D18").Value = 50
